$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 105; $r++) {
    $ws.Cells.Item($r, 3).Value = 7534
}

for ($r = 106; $r -le 124; $r++) {
    $ws.Cells.Item($r, 3).Value = 7320
}

for ($r = 125; $r -le 135; $r++) {
    $ws.Cells.Item($r, 3).Value = 7312
}

for ($r = 180; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
